$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A & C: strip trailing ".0" from numeric-looking ID/SKU strings ---
# Force text format so Excel COM keeps these as text (matching the source
# workbook's inline-string cells) instead of auto-converting to numbers.
$idFixes = @{
    "A2" = "3503"
    "A3" = "1845"
    "A4" = "3290"
    "A5" = "3289"
    "A6" = "3435"
    "A7" = "3434"
    "A8" = "5637"
    "A9" = "3504"
}
foreach ($addr in $idFixes.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $idFixes[$addr]
}

$skuFixes = @{
    "C2"  = "7798138551558"
    "C3"  = "7792180001665"
    "C4"  = "7797470133620"
    "C5"  = "7797470133637"
    "C6"  = "7798138551572"
    "C7"  = "7797470132739"
    "C8"  = "7792900806037"
    "C9"  = "7798138551589"
    "C10" = "7793065000117"
    "C11" = "7790380013525"
    "C12" = "7790380012061"
    "C13" = "7790380012122"
    "C14" = "7790380013068"
    "C15" = "7790380026044"
}
foreach ($addr in $skuFixes.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $skuFixes[$addr]
}

# --- Column F: refreshed "Nuevo Precio" values (auto-loaded Session IDs) ---
# Full literal text (including the leading "$" and non-breaking space) is
# stored directly so nothing gets auto-parsed as a currency number.
$priceFixes = @{
    "F4" = '$ 520,57'
    "F6" = '$ 578,43'
    "F7" = '$ 520,57'
    "F8" = '$ 1.652,81'
    "F9" = '$ 925,53'
}
foreach ($addr in $priceFixes.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $priceFixes[$addr]
}
